$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("COVERAGES")

# Delete row 2 (S_1, LU_1, 7.58) so the remaining data row (S_1, LU_1, 24.22)
# shifts up from row 3 to row 2, matching the updated default data set.
$ws.Rows.Item(2).Delete()
